# Daily attendance processing - 2025-12-25 21:53:55
# Swap the first two comma-separated entries in the "Recorded By" column (G)
# for every row that lists more than one recorder (any trailing entries
# keep their relative position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.Contains(",")) {
            $parts = $text -split ", "
            if ($parts.Count -ge 2) {
                $tmp = $parts[0]
                $parts[0] = $parts[1]
                $parts[1] = $tmp
                $cell.Value = $parts -join ", "
            }
        }
    }
}
